$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

function Append-RunToParagraph($paraIndex, $text) {
    $endPos = $d.Paragraphs($paraIndex).Range.End - 1
    $ins = $d.Range($endPos, $endPos)
    $ins.InsertAfter($text)
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-InParagraph 1 "ContosoLearn Competitor SWOT" "SWOT Kompetitor ContosoLearn"

# ---------------------------------------------------------------------------
# Fabrikam Learning section
# ---------------------------------------------------------------------------
Replace-InParagraph 2 "Fabrikam Learning:" "Pembelajaran Fabrikam:"

# Strengths
Replace-InParagraph 3 "Strengths:" "Kekuatan:"
Replace-InParagraph 3 " Fabrikam Learning provides a comprehensive set of analytics and reporting tools. It ensures the continuous monitoring of teaching and learning activities, as well as pinpointing problematic areas that need to be addressed." " Pembelajaran Fabrikam menyediakan serangkaian alat analitik dan pelaporan yang komprehensif."
Append-RunToParagraph 3 " "
Append-RunToParagraph 3 "Ini memastikan pemantauan berkelanjutan terhadap kegiatan belajar mengajar, serta menentukan area bermasalah yang perlu ditangani."

# Weaknesses
Replace-InParagraph 4 "Weaknesses:" "Kelemahan:"
Replace-InParagraph 4 " While Fabrikam Learning has robust reporting capabilities, it might be overwhelming for some users due to its comprehensive nature." " Meskipun Pembelajaran Fabrikam memiliki kemampuan pelaporan yang kuat, beberapa pengguna mungkin akan kewalahan karena sifatnya yang komprehensif."

# Opportunities
Replace-InParagraph 5 "Opportunities:" "Peluang:"
Replace-InParagraph 5 " There is a growing demand for personalized learning experiences and data-driven recommendations. Fabrikam Learning can leverage its robust analytics and reporting tools to meet this demand." " Ada peningkatan permintaan untuk pengalaman pembelajaran yang dipersonalisasi dan rekomendasi berbasis data."
Append-RunToParagraph 5 " "
Append-RunToParagraph 5 "Pembelajaran Fabrikam dapat memanfaatkan analitik dan alat pelaporan yang kuat untuk memenuhi permintaan ini."

# Threats
Replace-InParagraph 6 "Threats:" "Ancaman:"
Replace-InParagraph 6 " The eLearning market is highly competitive with many players offering similar features. Fabrikam Learning needs to continuously innovate to stay ahead." " Pasar eLearning sangat kompetitif dengan banyaknya pemain yang menawarkan fitur serupa."
Append-RunToParagraph 6 " "
Append-RunToParagraph 6 "Pembelajaran Fabrikam perlu terus berinovasi untuk tetap unggul."

# ---------------------------------------------------------------------------
# AdatumLearn section (heading text unchanged)
# ---------------------------------------------------------------------------

# Strengths
Replace-InParagraph 8 "Strengths:" "Kekuatan:"
Replace-InParagraph 8 " AdatumLearn offers courses on business analysis techniques such as MOST and SWOT. This shows their commitment to providing valuable content to their users." " AdatumLearn menawarkan kursus tentang teknik analisis bisnis seperti MOST dan SWOT."
Append-RunToParagraph 8 " "
Append-RunToParagraph 8 "Ini menunjukkan komitmen mereka untuk memberikan konten yang berharga bagi penggunanya."

# Weaknesses
Replace-InParagraph 9 "Weaknesses:" "Kelemahan:"
Replace-InParagraph 9 " The information provided in their courses is a compilation of third-party generated information. This might not be as valuable as original content." " Informasi yang diberikan dalam kursus mereka merupakan kompilasi informasi yang dihasilkan pihak ketiga."
Append-RunToParagraph 9 " "
Append-RunToParagraph 9 "Ini mungkin tidak seberharga konten asli."

# Opportunities
Replace-InParagraph 10 "Opportunities:" "Peluang:"
Replace-InParagraph 10 " AdatumLearn can create more original content to provide unique value to their users. They can also expand their course offerings to cover more topics." " AdatumLearn dapat membuat lebih banyak konten asli untuk memberikan nilai unik bagi penggunanya."
Append-RunToParagraph 10 " "
Append-RunToParagraph 10 "Mereka juga dapat memperluas penawaran kursusnya untuk mencakup lebih banyak topik."

# Threats
Replace-InParagraph 11 "Threats:" "Ancaman:"
Replace-InParagraph 11 " Like Fabrikam Learning, AdatumLearn also faces stiff competition in the eLearning market. They need to continuously improve their offerings to stay competitive.""" " Seperti Pembelajaran Fabrikam, AdatumLearn juga menghadapi persaingan ketat di pasar eLearning."
Append-RunToParagraph 11 " "
Append-RunToParagraph 11 "Mereka perlu meningkatkan penawarannya secara terus-menerus agar tetap kompetitif."""
